# Auto-generated: applies the per-cell Price (D) / Volume(1h) (E) text updates
# from the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text (e.g. "23.420.47", "303.33") in the source data,
# not numbers -- force text formatting first so Excel does not auto-convert
# numeric-looking strings (like "303.33") into floating point numbers and lose
# the original formatting (e.g. "1.000" must stay "1.000", not become 1).
$dCells = @('D2', 'D3', 'D6', 'D7', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '23.420.47'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '1.634.59'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = '303.33'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('D7').Value = '0.3781'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = '51.59'
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('D10').Value = '0.08179'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '1.227'
$ws.Range('E11').Value = '  -3.67%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '22.34'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '6.460'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').Value = '7.376'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '0.00001238'
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').Value = '1.629.84'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '95.03'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '0.06941'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').Value = '6.582'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '17.44'
$ws.Range('E21').Value = '  -5.45%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('D24').Value = '23.429.01'
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('D25').Value = '2.517'
$ws.Range('D26').Value = '3.044'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('D27').Value = '21.12'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').Value = '150.88'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '5.270'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').Value = '133.39'
$ws.Range('E30').Value = '  -2.87%  '
$ws.Range('D31').Value = '1.812.34'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '2.167'
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('D33').Value = '6.598'
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('D34').Value = '1.045'
$ws.Range('D35').Value = '11.21'
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('D36').Value = '0.02745'
$ws.Range('E36').Value = '  -4.53%  '
$ws.Range('D37').Value = '0.08768'
$ws.Range('D38').Value = '0.2491'
$ws.Range('E38').Value = '  -3.33%  '
$ws.Range('D39').Value = '0.07104'
$ws.Range('E39').Value = '  -3.85%  '
$ws.Range('D40').Value = '6.011'
$ws.Range('E40').Value = '  -5.86%  '
$ws.Range('D41').Value = '0.7005'
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('D42').Value = '1.339'
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').Value = '15.82'
$ws.Range('E43').Value = '  -4.29%  '
$ws.Range('D44').Value = '12.13'
$ws.Range('E44').Value = '  -4.75%  '
$ws.Range('D45').Value = '0.6491'
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').Value = '2.271'
$ws.Range('E47').Value = '  -4.45%  '
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('D49').Value = '0.07973'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').Value = '127.12'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '1.189'
$ws.Range('E51').Value = '  -3.69%  '
